$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 45947
$ws.Range("B2").Value = 6096.97941642731
$ws.Range("C2").Value = 5617.01037104988
$ws.Range("D2").Value = 8852
$ws.Range("E2").Value = 8773.036039
$ws.Range("F2").Value = -23.2888752657263

# Row 3
$ws.Range("A3").Value = 45948
$ws.Range("B3").Value = 2066.82194182061
$ws.Range("C3").Value = 4138.43869325176
$ws.Range("D3").Value = 3620
$ws.Range("E3").Value = 4744.499386
$ws.Range("F3").Value = 133.171505726298

# Row 4
$ws.Range("A4").Value = 45949
$ws.Range("B4").Value = 2048.56673871219
$ws.Range("C4").Value = 4204.32527585282
$ws.Range("D4").Value = 3620
$ws.Range("E4").Value = 4802.039946
$ws.Range("F4").Value = 139.074936797526

# Row 5
$ws.Range("A5").Value = 45950
$ws.Range("B5").Value = 7245.65739014288
$ws.Range("C5").Value = 6983.25384047283
$ws.Range("D5").Value = 3620
$ws.Range("E5").Value = 10975.050045
$ws.Range("F5").Value = 295.526937305414

# Row 6
$ws.Range("A6").Value = 45951
$ws.Range("B6").Value = 7696.89866058229
$ws.Range("C6").Value = 7125.71832999126
$ws.Range("D6").Value = 3620
$ws.Range("E6").Value = 11781.04191
$ws.Range("F6").Value = 316.244232475374

# Row 7
$ws.Range("A7").Value = 45952
$ws.Range("B7").Value = 6392.95297294923
$ws.Range("C7").Value = 6381.84361816817
$ws.Range("D7").Value = 3620
$ws.Range("E7").Value = 9928.531221
$ws.Range("F7").Value = 262.392577759122

# Row 8
$ws.Range("A8").Value = 45953
$ws.Range("B8").Value = 6392.95297294923
$ws.Range("C8").Value = 6080.27928711116
$ws.Range("D8").Value = 3620
$ws.Range("E8").Value = 9928.531221
$ws.Range("F8").Value = 249.827397298414

# Row 9
$ws.Range("A9").Value = 45954
$ws.Range("B9").Value = 6392.95297294923
$ws.Range("C9").Value = 5260.05113413809
$ws.Range("D9").Value = 3620
$ws.Range("E9").Value = 9928.531221
$ws.Range("F9").Value = 215.651224257869

# Row 10
$ws.Range("A10").Value = 45955
$ws.Range("B10").Value = 2044.14689306201
$ws.Range("C10").Value = 3842.36053325781
$ws.Range("D10").Value = 3620
$ws.Range("E10").Value = 5157.791021
$ws.Range("F10").Value = 139.000194216491

# Row 11
$ws.Range("A11").Value = 45956
$ws.Range("B11").Value = 1933.9611505005
$ws.Range("C11").Value = 3730.16885780957
$ws.Range("D11").Value = 3620
$ws.Range("E11").Value = 5037.344132
$ws.Range("F11").Value = 133.897993304545

# Row 12
$ws.Range("A12").Value = 45957
$ws.Range("B12").Value = 6757.64823424925
$ws.Range("C12").Value = 6546.80289865024
$ws.Range("D12").Value = 3620
$ws.Range("E12").Value = 10717.736825
$ws.Range("F12").Value = 286.953812058375

# Row 13
$ws.Range("A13").Value = 45958
$ws.Range("B13").Value = 6757.64823424925
$ws.Range("C13").Value = 7042.0933831584
$ws.Range("D13").Value = 3620
$ws.Range("E13").Value = 10717.736825
$ws.Range("F13").Value = 307.590915579548

# Row 14
$ws.Range("A14").Value = 45959
$ws.Range("B14").Value = 6757.64823424925
$ws.Range("C14").Value = 7224.87941415647
$ws.Range("D14").Value = 3620
$ws.Range("E14").Value = 10717.736825
$ws.Range("F14").Value = 315.207000204468

# Row 15
$ws.Range("A15").Value = 45960
$ws.Range("B15").Value = 6757.64823424925
$ws.Range("C15").Value = 7378.44428301094
$ws.Range("D15").Value = 3620
$ws.Range("E15").Value = 10717.736825
$ws.Range("F15").Value = 321.605536406737

Write-Output "done"
